$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")
$ws.Activate()

$ws.Range("A4").Value = "LOG030"
$ws.Range("B4").Value = "log cancel"

$ws.Range("A5").Select()
